$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.523673
$ws.Range("H2").Value = 1.571019
$ws.Range("I2").Value = 0.001411687926676084
$ws.Range("J2").Value = 0.001412631582657676
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 109.1447706666667
$ws.Range("N2").Value = 327.434312
$ws.Range("O2").Value = 0.3535542089399963
$ws.Range("P2").Value = 0.3655959674582361
$ws.Range("Q2").Value = 57.15616948932534
$ws.Range("R2").Value = 514.405525403928
$ws.Range("S2").Value = 0.0004991082081861062
$ws.Range("T2").Value = 0.0005164524101237923

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.523673
$ws.Range("H3").Value = 1.571019
$ws.Range("I3").Value = 0.001411687926676084
$ws.Range("J3").Value = 0.001412631582657676
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 47.980825
$ws.Range("N3").Value = 143.942475
$ws.Range("O3").Value = 0.155424969272891
$ws.Range("P3").Value = 0.1607186127944892
$ws.Range("Q3").Value = 25.12626257022501
$ws.Range("R3").Value = 226.136363132025
$ws.Range("S3").Value = 0.0002194115526265415
$ws.Range("T3").Value = 0.0002270361883544255

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.523673
$ws.Range("H4").Value = 1.571019
$ws.Range("I4").Value = 0.001411687926676084
$ws.Range("J4").Value = 0.001412631582657676
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.26741999999999
$ws.Range("N4").Value = 186.80226
$ws.Range("O4").Value = 0.2017037397794264
$ws.Range("P4").Value = 0.2085735992386923
$ws.Range("Q4").Value = 32.60776663366
$ws.Range("R4").Value = 293.46989970294
$ws.Range("S4").Value = 0.0002847427342120307
$ws.Range("T4").Value = 0.0002946376535931618

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.523673
$ws.Range("H5").Value = 1.571019
$ws.Range("I5").Value = 0.001411687926676084
$ws.Range("J5").Value = 0.001412631582657676
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 58.81030666666667
$ws.Range("N5").Value = 176.43092
$ws.Range("O5").Value = 0.1905050633580386
$ws.Range("P5").Value = 0.1969935053322898
$ws.Range("Q5").Value = 30.79736972305334
$ws.Range("R5").Value = 277.17632750748
$ws.Range("S5").Value = 0.0002689336979132054
$ws.Range("T5").Value = 0.000278279247210836

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.523673
$ws.Range("H6").Value = 1.571019
$ws.Range("I6").Value = 0.001411687926676084
$ws.Range("J6").Value = 0.001412631582657676
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.503993
$ws.Range("N6").Value = 61.007986
$ws.Range("O6").Value = 0.09881201864964768
$ws.Range("P6").Value = 0.06811831517629259
$ws.Range("Q6").Value = 15.974117526289
$ws.Range("R6").Value = 95.84470515773401
$ws.Range("S6").Value = 0.0001394917337381996
$ws.Range("T6").Value = 0.0000962260833754606

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 290.8976950000001
$ws.Range("H7").Value = 872.6930850000001
$ws.Range("I7").Value = 0.7841854820267643
$ws.Range("J7").Value = 0.78470967813754
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 109.1447706666667
$ws.Range("N7").Value = 327.434312
$ws.Range("O7").Value = 0.3535542089399963
$ws.Range("P7").Value = 0.3655959674582361
$ws.Range("Q7").Value = 31749.96220823695
$ws.Range("R7").Value = 285749.6598741325
$ws.Range("S7").Value = 0.2772520777602024
$ws.Range("T7").Value = 0.286886693952535

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 290.8976950000001
$ws.Range("H8").Value = 872.6930850000001
$ws.Range("I8").Value = 0.7841854820267643
$ws.Range("J8").Value = 0.78470967813754
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 47.980825
$ws.Range("N8").Value = 143.942475
$ws.Range("O8").Value = 0.155424969272891
$ws.Range("P8").Value = 0.1607186127944892
$ws.Range("Q8").Value = 13957.51139669838
$ws.Range("R8").Value = 125617.6025702854
$ws.Range("S8").Value = 0.121882004448257
$ws.Range("T8").Value = 0.1261174509166755

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 290.8976950000001
$ws.Range("H9").Value = 872.6930850000001
$ws.Range("I9").Value = 0.7841854820267643
$ws.Range("J9").Value = 0.78470967813754
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 62.26741999999999
$ws.Range("N9").Value = 186.80226
$ws.Range("O9").Value = 0.2017037397794264
$ws.Range("P9").Value = 0.2085735992386923
$ws.Range("Q9").Value = 18113.4489515969
$ws.Range("R9").Value = 163021.0405643721
$ws.Range("S9").Value = 0.1581731444055305
$ws.Range("T9").Value = 0.1636697219265825

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 290.8976950000001
$ws.Range("H10").Value = 872.6930850000001
$ws.Range("I10").Value = 0.7841854820267643
$ws.Range("J10").Value = 0.78470967813754
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 58.81030666666667
$ws.Range("N10").Value = 176.43092
$ws.Range("O10").Value = 0.1905050633580386
$ws.Range("P10").Value = 0.1969935053322898
$ws.Range("Q10").Value = 17107.78265157647
$ws.Range("R10").Value = 153970.0438641882
$ws.Range("S10").Value = 0.1493913049379628
$ws.Range("T10").Value = 0.1545827101644869

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 290.8976950000001
$ws.Range("H11").Value = 872.6930850000001
$ws.Range("I11").Value = 0.7841854820267643
$ws.Range("J11").Value = 0.78470967813754
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 30.503993
$ws.Range("N11").Value = 61.007986
$ws.Range("O11").Value = 0.09881201864964768
$ws.Range("P11").Value = 0.06811831517629259
$ws.Range("Q11").Value = 8873.541251996137
$ws.Range("R11").Value = 53241.24751197682
$ws.Range("S11").Value = 0.07748695047481159
$ws.Range("T11").Value = 0.05345310117726006

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 38.68424166666667
$ws.Range("H12").Value = 116.052725
$ws.Range("I12").Value = 0.1042827812651277
$ws.Range("J12").Value = 0.1043524900644015
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 109.1447706666667
$ws.Range("N12").Value = 327.434312
$ws.Range("O12").Value = 0.3535542089399963
$ws.Range("P12").Value = 0.3655959674582361
$ws.Range("Q12").Value = 4222.182685122244
$ws.Range("R12").Value = 37999.6441661002
$ws.Range("S12").Value = 0.0368696162362549
$ws.Range("T12").Value = 0.03815084956177085

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 38.68424166666667
$ws.Range("H13").Value = 116.052725
$ws.Range("I13").Value = 0.1042827812651277
$ws.Range("J13").Value = 0.1043524900644015
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 47.980825
$ws.Range("N13").Value = 143.942475
$ws.Range("O13").Value = 0.155424969272891
$ws.Range("P13").Value = 0.1607186127944892
$ws.Range("Q13").Value = 1856.101829666042
$ws.Range("R13").Value = 16704.91646699437
$ws.Range("S13").Value = 0.01620814807382408
$ws.Range("T13").Value = 0.01677138744480133

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 38.68424166666667
$ws.Range("H14").Value = 116.052725
$ws.Range("I14").Value = 0.1042827812651277
$ws.Range("J14").Value = 0.1043524900644015
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 62.26741999999999
$ws.Range("N14").Value = 186.80226
$ws.Range("O14").Value = 0.2017037397794264
$ws.Range("P14").Value = 0.2085735992386923
$ws.Range("Q14").Value = 2408.767923239833
$ws.Range("R14").Value = 21678.9113091585
$ws.Range("S14").Value = 0.02103422697577616
$ws.Range("T14").Value = 0.02176517444225211

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 38.68424166666667
$ws.Range("H15").Value = 116.052725
$ws.Range("I15").Value = 0.1042827812651277
$ws.Range("J15").Value = 0.1043524900644015
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 58.81030666666667
$ws.Range("N15").Value = 176.43092
$ws.Range("O15").Value = 0.1905050633580386
$ws.Range("P15").Value = 0.1969935053322898
$ws.Range("Q15").Value = 2275.032115584111
$ws.Range("R15").Value = 20475.289040257
$ws.Range("S15").Value = 0.01986639785206563
$ws.Range("T15").Value = 0.02055676280793941

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 38.68424166666667
$ws.Range("H16").Value = 116.052725
$ws.Range("I16").Value = 0.1042827812651277
$ws.Range("J16").Value = 0.1043524900644015
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 30.503993
$ws.Range("N16").Value = 61.007986
$ws.Range("O16").Value = 0.09881201864964768
$ws.Range("P16").Value = 0.06811831517629259
$ws.Range("Q16").Value = 1180.023837010308
$ws.Range("R16").Value = 7080.14302206185
$ws.Range("S16").Value = 0.01030439212720693
$ws.Range("T16").Value = 0.007108315807637844

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 40.106198
$ws.Range("H17").Value = 120.318594
$ws.Range("I17").Value = 0.1081160103757125
$ws.Range("J17").Value = 0.1081882815327926
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 109.1447706666667
$ws.Range("N17").Value = 327.434312
$ws.Range("O17").Value = 0.3535542089399963
$ws.Range("P17").Value = 0.3655959674582361
$ws.Range("Q17").Value = 4377.381783021925
$ws.Range("R17").Value = 39396.43604719732
$ws.Range("S17").Value = 0.03822487052213346
$ws.Range("T17").Value = 0.03955319945462533

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 40.106198
$ws.Range("H18").Value = 120.318594
$ws.Range("I18").Value = 0.1081160103757125
$ws.Range("J18").Value = 0.1081882815327926
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 47.980825
$ws.Range("N18").Value = 143.942475
$ws.Range("O18").Value = 0.155424969272891
$ws.Range("P18").Value = 0.1607186127944892
$ws.Range("Q18").Value = 1924.32846765335
$ws.Range("R18").Value = 17318.95620888015
$ws.Range("S18").Value = 0.01680392759055267
$ws.Range("T18").Value = 0.01738787052857008

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 40.106198
$ws.Range("H19").Value = 120.318594
$ws.Range("I19").Value = 0.1081160103757125
$ws.Range("J19").Value = 0.1081882815327926
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 62.26741999999999
$ws.Range("N19").Value = 186.80226
$ws.Range("O19").Value = 0.2017037397794264
$ws.Range("P19").Value = 0.2085735992386923
$ws.Range("Q19").Value = 2497.30947546916
$ws.Range("R19").Value = 22475.78527922244
$ws.Range("S19").Value = 0.02180740362281247
$ws.Range("T19").Value = 0.02256521927474351

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 40.106198
$ws.Range("H20").Value = 120.318594
$ws.Range("I20").Value = 0.1081160103757125
$ws.Range("J20").Value = 0.1081882815327926
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 58.81030666666667
$ws.Range("N20").Value = 176.43092
$ws.Range("O20").Value = 0.1905050633580386
$ws.Range("P20").Value = 0.1969935053322898
$ws.Range("Q20").Value = 2358.657803614054
$ws.Range("R20").Value = 21227.92023252648
$ws.Range("S20").Value = 0.02059664740664346
$ws.Range("T20").Value = 0.02131238881502146

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 40.106198
$ws.Range("H21").Value = 120.318594
$ws.Range("I21").Value = 0.1081160103757125
$ws.Range("J21").Value = 0.1081882815327926
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 30.503993
$ws.Range("N21").Value = 61.007986
$ws.Range("O21").Value = 0.09881201864964768
$ws.Range("P21").Value = 0.06811831517629259
$ws.Range("Q21").Value = 1223.399183048614
$ws.Range("R21").Value = 7340.395098291684
$ws.Range("S21").Value = 0.0106831612335704
$ws.Range("T21").Value = 0.007369603459832243

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 0.7434085
$ws.Range("H22").Value = 1.486817
$ws.Range("I22").Value = 0.002004038405719556
$ws.Range("J22").Value = 0.001336918682608127
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 109.1447706666667
$ws.Range("N22").Value = 327.434312
$ws.Range("O22").Value = 0.3535542089399963
$ws.Range("P22").Value = 0.3655959674582361
$ws.Range("Q22").Value = 81.13915024415067
$ws.Range("R22").Value = 486.834901464904
$ws.Range("S22").Value = 0.0007085362132195491
$ws.Range("T22").Value = 0.0004887720791811088

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 0.7434085
$ws.Range("H23").Value = 1.486817
$ws.Range("I23").Value = 0.002004038405719556
$ws.Range("J23").Value = 0.001336918682608127
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 47.980825
$ws.Range("N23").Value = 143.942475
$ws.Range("O23").Value = 0.155424969272891
$ws.Range("P23").Value = 0.1607186127944892
$ws.Range("Q23").Value = 35.6693531420125
$ws.Range("R23").Value = 214.016118852075
$ws.Range("S23").Value = 0.0003114776076306554
$ws.Range("T23").Value = 0.0002148677160878142

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 0.7434085
$ws.Range("H24").Value = 1.486817
$ws.Range("I24").Value = 0.002004038405719556
$ws.Range("J24").Value = 0.001336918682608127
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 62.26741999999999
$ws.Range("N24").Value = 186.80226
$ws.Range("O24").Value = 0.2017037397794264
$ws.Range("P24").Value = 0.2085735992386923
$ws.Range("Q24").Value = 46.29012930107
$ws.Range("R24").Value = 277.74077580642
$ws.Range("S24").Value = 0.0004042220410952339
$ws.Range("T24").Value = 0.0002788459415210281

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 0.7434085
$ws.Range("H25").Value = 1.486817
$ws.Range("I25").Value = 0.002004038405719556
$ws.Range("J25").Value = 0.001336918682608127
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 58.81030666666667
$ws.Range("N25").Value = 176.43092
$ws.Range("O25").Value = 0.1905050633580386
$ws.Range("P25").Value = 0.1969935053322898
$ws.Range("Q25").Value = 43.72008186360667
$ws.Range("R25").Value = 262.32049118164
$ws.Range("S25").Value = 0.0003817794634535467
$ws.Range("T25").Value = 0.0002633642976312021

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 0.7434085
$ws.Range("H26").Value = 1.486817
$ws.Range("I26").Value = 0.002004038405719556
$ws.Range("J26").Value = 0.001336918682608127
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 30.503993
$ws.Range("N26").Value = 61.007986
$ws.Range("O26").Value = 0.09881201864964768
$ws.Range("P26").Value = 0.06811831517629259
$ws.Range("Q26").Value = 22.6769276801405
$ws.Range("R26").Value = 90.707710720562
$ws.Range("S26").Value = 0.000198023080320571
$ws.Range("T26").Value = 0.0000910686481869743
